$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 618 (shifts rows 618:639 down to 620:641),
# mirroring the weekly update that prepends the newest price observations.
$ws.Rows.Item(618).Insert()
$ws.Rows.Item(618).Insert()

# New row 618: Femacal de La Calera, Limon "1a amarillo", week of 2021-11-09
$ws.Cells.Item(618, 1).Value = 3
$ws.Cells.Item(618, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(618, 3).Value = "Coquimbo"
$ws.Cells.Item(618, 4).Value = 44509
$ws.Cells.Item(618, 5).Value = 5
$ws.Cells.Item(618, 6).Value = "Fruta"
$ws.Cells.Item(618, 7).Value = 100102
$ws.Cells.Item(618, 8).Value = "Cítricos"
$ws.Cells.Item(618, 9).Value = 100102003
$ws.Cells.Item(618, 10).Value = "Limón"
$ws.Cells.Item(618, 11).Value = "Sin especificar"
$ws.Cells.Item(618, 12).Value = "1a amarillo"
$ws.Cells.Item(618, 13).Value = 212
$ws.Cells.Item(618, 14).Value = 4500
$ws.Cells.Item(618, 15).Value = 5000
$ws.Cells.Item(618, 16).Value = 4764
$ws.Cells.Item(618, 17).Value = "$/malla 16 kilos"
$ws.Cells.Item(618, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(618, 19).Value = 298
$ws.Cells.Item(618, 20).Value = 16

# New row 619: Femacal de La Calera, Limon "2a amarillo", week of 2021-11-09
$ws.Cells.Item(619, 1).Value = 3
$ws.Cells.Item(619, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(619, 3).Value = "Coquimbo"
$ws.Cells.Item(619, 4).Value = 44509
$ws.Cells.Item(619, 5).Value = 5
$ws.Cells.Item(619, 6).Value = "Fruta"
$ws.Cells.Item(619, 7).Value = 100102
$ws.Cells.Item(619, 8).Value = "Cítricos"
$ws.Cells.Item(619, 9).Value = 100102003
$ws.Cells.Item(619, 10).Value = "Limón"
$ws.Cells.Item(619, 11).Value = "Sin especificar"
$ws.Cells.Item(619, 12).Value = "2a amarillo"
$ws.Cells.Item(619, 13).Value = 210
$ws.Cells.Item(619, 14).Value = 3500
$ws.Cells.Item(619, 15).Value = 4000
$ws.Cells.Item(619, 16).Value = 3762
$ws.Cells.Item(619, 17).Value = "$/malla 16 kilos"
$ws.Cells.Item(619, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(619, 19).Value = 235
$ws.Cells.Item(619, 20).Value = 16
